$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entregables")

# Mark row 2 (agregarArista...) as done
$ws.Range("C2").Value = "x"

# Assign responsible persons in column D
$ws.Range("D2").Value = "Javier"
$ws.Range("D9").Value = "Camila"
$ws.Range("D10").Value = "Camila"

# Move selection to D12 as last active cell
$ws.Range("D12").Select()
